# Remove unused dt and dx parameters in favour of the higher level number of
# timesteps given from the excel parameter sheet (general!normal_timestep /
# general!shorter_timestep), and add explicit "dx" grid-spacing parameters on
# the channel and peat sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "general" sheet: add normal_timestep / shorter_timestep columns
# ---------------------------------------------------------------------------
$general = $wb.Worksheets.Item("general")

$general.Range("B1").Value = "normal_timestep"
$general.Range("C1").Value = "shorter_timestep"
$general.Range("B2").Value = 24
$general.Range("C2").Value = 1000

$bComment = $general.Range("B1").AddComment("Dt = 1day/timestep;`nI.e., number of timesteps per day. For instance, a value of 24 means that dt=1 hour.`n ")
$cComment = $general.Range("C1").AddComment("Used as alternative timestep when normal_timestep doesnt converge")

# ---------------------------------------------------------------------------
# "channel" sheet: insert a "dx" column between channel_width and
# porous_threshold, and highlight the last two (max_niter_newton*) columns
# ---------------------------------------------------------------------------
$channel = $wb.Worksheets.Item("channel")

$channel.Columns.Item(9).Insert()
$channel.Range("I1").Value = "dx"
$channel.Range("I2").Value = 100

$channel.Range("N1:O2").Interior.Color = 15658734

# ---------------------------------------------------------------------------
# "peat" sheet: add a "dx" column (grid spacing, metres)
# ---------------------------------------------------------------------------
$peat = $wb.Worksheets.Item("peat")

$peat.Range("D1").Value = "dx"
$peat.Range("D2").Value = 50

$dComment = $peat.Range("D1").AddComment("Metres. Only used if rectangular grid`n")

# ---------------------------------------------------------------------------
# Selections / active sheet bookkeeping
# ---------------------------------------------------------------------------
$channel.Range("I4").Select()

$peat.Range("D1").Select()

$general.Activate()
$general.Range("C1").Select()
